# "Implementação básica do login." — add a new team member (Hygor) and a
# "Cancelado" status, insert the new front-end tasks he's responsible for,
# and refresh the status of a few in-flight tasks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected (no password) - unprotect for the duration of
# the edit and re-protect once we're done.
$ws.Unprotect()

# --- 1. Extend the two helper lists used by the data-validation dropdowns ---
# Names list (C column source) gains "Hygor" right after "Rodrigo", and the
# status list (F column source) gains "Cancelado" after "Concluído".
$ws.Range("K8").Value = "Hygor"
$ws.Range("K15").Value = "Cancelado"

# --- 2. Insert 9 new rows for the front-end tasks, right before the old
#        "Usuário CRUD" row (old row 16), pushing everything below down. ---
$ws.Range("A15:F15").Copy()
$ws.Rows.Item(16).Resize(9).Insert(-4121)
$ws.Range("A16:F24").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A16").Value = "Front-end Tela de login"
$ws.Range("C16").Value = "Luiz Henrique"
$ws.Range("D16").Value = 0
$ws.Range("F16").Value = "Pausado"

$ws.Range("A17").Value = "Front-end Tela de cadastro/alteração de usuário"
$ws.Range("C17").Value = "Luiz Henrique"
$ws.Range("D17").Value = 0
$ws.Range("F17").Value = "Pendente"

$ws.Range("A18").Value = "Front-end Tela de recuperação de senha"
$ws.Range("C18").Value = "Luiz Henrique"
$ws.Range("D18").Value = 0
$ws.Range("F18").Value = "Pendente"

$ws.Range("A19").Value = "Front-end Tela de cadastro/alteração de veículo"
$ws.Range("C19").Value = "Luiz Henrique"
$ws.Range("D19").Value = 0
$ws.Range("F19").Value = "Pendente"

$ws.Range("A20").Value = "Front-end Tela de cadastro/alteração de nota fiscal"
$ws.Range("C20").Value = "Luiz Henrique"
$ws.Range("D20").Value = 0
$ws.Range("F20").Value = "Pendente"

$ws.Range("A21").Value = "Front-end Tela de listagem/exclusão das notas fiscais"
$ws.Range("C21").Value = "Luiz Henrique"
$ws.Range("D21").Value = 0
$ws.Range("F21").Value = "Pendente"

$ws.Range("A22").Value = "Front-end Tela de gerenciamento de roteiro"
$ws.Range("C22").Value = "Luiz Henrique"
$ws.Range("D22").Value = 0
$ws.Range("F22").Value = "Pendente"

$ws.Range("A23").Value = "Front-end Tela do mapa"
$ws.Range("C23").Value = "Luiz Henrique"
$ws.Range("D23").Value = 0
$ws.Range("F23").Value = "Pendente"

$ws.Range("A24").Value = "Front-end Tela de visualização dos roteiros (para o gerente)"
$ws.Range("C24").Value = "Luiz Henrique"
$ws.Range("D24").Value = 0
$ws.Range("F24").Value = "Pendente"

# The "Usuário CRUD" row (now row 25) gained an executante too.
$ws.Range("C25").Value = "Luiz Henrique"

# --- 3. Refresh statuses on a few existing tasks. ---
$ws.Range("F5").Value = "Cancelado"
$ws.Range("F7").Value = "Concluído"
$ws.Range("F8").Value = "Em andamento"
$ws.Range("F13").Value = "Em andamento"
$ws.Range("F25").Value = "Concluído"

# --- 4. Data validation sources need to track the widened helper lists. ---
$ws.Range("F3:F32").Validation.Delete()
$ws.Range("F3:F32").Validation.Add(3, 1, 1, "=`$K`$10:`$K`$15")

$ws.Range("C3:C32").Validation.Delete()
$ws.Range("C3:C32").Validation.Add(3, 1, 1, "=`$K`$3:`$K`$8")

# --- 6. Conditional formatting: add a "Cancelado" rule, highest priority,
#        matching the style used for the other status rules (white bold
#        text over a themed dark fill). ---
$ccf = $ws.Range("F3:F32").FormatConditions
$newRule = $ccf.Add(1, 3, """Cancelado""")
$newRule.SetFirstPriority()
$newRule.Font.Bold = $true
$newRule.Font.Italic = $false
$newRule.Font.ColorIndex = 1
$newRule.Interior.ThemeColor = 1
$newRule.Interior.TintAndShade = -0.499984740745262

$ws.Range("F31").Select()

$ws.Protect()
